$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 26
$ws.Range("D2").Value = 44657
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1857

# Row 3 <- source row 11
$ws.Range("D3").Value = 44623
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 2286

# Row 4 <- source row 22
$ws.Range("D4").Value = 44690
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1714

# Row 5 <- source row 32
$ws.Range("D5").Value = 44302
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 340
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("R5").Value = "Provincia de Santiago"
$ws.Range("S5").Value = 1786

# Row 6 <- source row 18
$ws.Range("D6").Value = 44312
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 1857

# Row 7 <- source row 19
$ws.Range("D7").Value = 44312
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 11000
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1571

# Row 8 <- source row 7
$ws.Range("D8").Value = 44659
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 2143

# Row 9 <- source row 8
$ws.Range("D9").Value = 44659
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1714

# Row 10 <- source row 24
$ws.Range("D10").Value = 44306
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 1714

# Row 11 <- source row 25
$ws.Range("D11").Value = 44306
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1286

# Row 12 <- source row 14
$ws.Range("D12").Value = 44687
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 2143

# Row 13 <- source row 15
$ws.Range("D13").Value = 44687
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 75
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 1714

# Row 14 <- source row 20
$ws.Range("D14").Value = 44307
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 2000

# Row 15 <- source row 21
$ws.Range("D15").Value = 44307
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 1429

# Row 16 <- source row 36
$ws.Range("D16").Value = 44322
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 11000
$ws.Range("R16").Value = "Región Metropolitana"
$ws.Range("S16").Value = 1571

# Row 17 <- source row 9
$ws.Range("D17").Value = 44300
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 12500
$ws.Range("R17").Value = "Provincia de Santiago"
$ws.Range("S17").Value = 1786

# Row 18 <- source row 16
$ws.Range("D18").Value = 44321
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = 11000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 11500
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1643

# Row 19 <- source row 17
$ws.Range("D19").Value = 44321
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("S19").Value = 1143

# Row 20 <- source row 40
$ws.Range("D20").Value = 44644
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 85
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("R20").Value = "Región Metropolitana"
$ws.Range("S20").Value = 2000

# Row 21 <- source row 39
$ws.Range("D21").Value = 44643
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("R21").Value = "Región Metropolitana"
$ws.Range("S21").Value = 2143

# Row 22 <- source row 2
$ws.Range("D22").Value = 44335
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("R22").Value = "Región Metropolitana"
$ws.Range("S22").Value = 2000

# Row 23 <- source row 12
$ws.Range("D23").Value = 44685
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 15000
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 2143

# Row 24 <- source row 13
$ws.Range("D24").Value = 44685
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 70
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("R24").Value = "Región Metropolitana"
$ws.Range("S24").Value = 1714

# Row 25 <- source row 27
$ws.Range("D25").Value = 44314
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 20
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 1857

# Row 26 <- source row 28
$ws.Range("D26").Value = 44314
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 45
$ws.Range("N26").Value = 11000
$ws.Range("O26").Value = 11000
$ws.Range("P26").Value = 11000
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 1571

# Row 27 <- source row 33
$ws.Range("D27").Value = 44349
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 70
$ws.Range("N27").Value = 12000
$ws.Range("O27").Value = 12000
$ws.Range("P27").Value = 12000
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 1714

# Row 28 <- source row 37
$ws.Range("D28").Value = 44694
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 80
$ws.Range("N28").Value = 15000
$ws.Range("O28").Value = 15000
$ws.Range("P28").Value = 15000
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 2143

# Row 29 <- source row 38
$ws.Range("D29").Value = 44694
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 75
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 1714

# Row 30 <- source row 34
$ws.Range("D30").Value = 44316
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = 13000
$ws.Range("O30").Value = 13000
$ws.Range("P30").Value = 13000
$ws.Range("R30").Value = "Región Metropolitana"
$ws.Range("S30").Value = 1857

# Row 31 <- source row 35
$ws.Range("D31").Value = 44316
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 11000
$ws.Range("O31").Value = 11000
$ws.Range("P31").Value = 11000
$ws.Range("R31").Value = "Región Metropolitana"
$ws.Range("S31").Value = 1571

# Row 32 <- source row 3
$ws.Range("D32").Value = 44342
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 12000
$ws.Range("O32").Value = 12000
$ws.Range("P32").Value = 12000
$ws.Range("R32").Value = "Región Metropolitana"
$ws.Range("S32").Value = 1714

# Row 33 <- source row 23
$ws.Range("D33").Value = 44641
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 13000
$ws.Range("O33").Value = 13000
$ws.Range("P33").Value = 13000
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 1857

# Row 34 <- source row 10
$ws.Range("D34").Value = 44679
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 150
$ws.Range("N34").Value = 12000
$ws.Range("O34").Value = 12000
$ws.Range("P34").Value = 12000
$ws.Range("R34").Value = "Región Metropolitana"
$ws.Range("S34").Value = 1714

# Row 35 <- source row 4
$ws.Range("D35").Value = 44664
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 80
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 14000
$ws.Range("P35").Value = 14000
$ws.Range("R35").Value = "Región Metropolitana"
$ws.Range("S35").Value = 2000

# Row 36 <- source row 5
$ws.Range("D36").Value = 44664
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 12000
$ws.Range("P36").Value = 12000
$ws.Range("R36").Value = "Región Metropolitana"
$ws.Range("S36").Value = 1714

# Row 37 <- source row 6
$ws.Range("D37").Value = 44344
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 12000
$ws.Range("R37").Value = "Región Metropolitana"
$ws.Range("S37").Value = 1714

# Row 38 <- source row 29
$ws.Range("D38").Value = 44315
$ws.Range("L38").Value = "Especial"
$ws.Range("M38").Value = 50
$ws.Range("N38").Value = 14000
$ws.Range("O38").Value = 14000
$ws.Range("P38").Value = 14000
$ws.Range("R38").Value = "Región Metropolitana"
$ws.Range("S38").Value = 2000

# Row 39 <- source row 30
$ws.Range("D39").Value = 44315
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 80
$ws.Range("N39").Value = 12000
$ws.Range("O39").Value = 13000
$ws.Range("P39").Value = 12500
$ws.Range("R39").Value = "Región Metropolitana"
$ws.Range("S39").Value = 1786

# Row 40 <- source row 31
$ws.Range("D40").Value = 44315
$ws.Range("L40").Value = "Segunda"
$ws.Range("M40").Value = 80
$ws.Range("N40").Value = 10000
$ws.Range("O40").Value = 11000
$ws.Range("P40").Value = 10500
$ws.Range("R40").Value = "Región Metropolitana"
$ws.Range("S40").Value = 1500
